$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting from the existing header cell (H1) onto the new
# header cells so they pick up the same style index as the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Header text for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-25
$data = @(
    @(8, 8),
    @(7, 8),
    @(8, 8),
    @(8, 8),
    @(6, 6),
    @(8, 9),
    @(6, 7),
    @(8, 8),
    @(6, 7),
    @(7, 8),
    @(7, 7),
    @(5, 5),
    @(8, 8),
    @(8, 8),
    @(7, 7),
    @(8, 8),
    @(6, 6),
    @(5, 5),
    @(5, 5),
    @(1, 1),
    @(4, 5),
    @(6, 6),
    @(7, 7),
    @(9, 9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $data[$i][0]
    $ws.Cells.Item($row, 10).Value = $data[$i][1]
}
